$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 712, shifting the existing rows 712:735 down to 713:736.
$ws.Rows(712).Insert()

# Populate the newly inserted row 712 with its data.
$ws.Range("A712").Value = 6
$ws.Range("B712").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C712").Value = "Metropolitana"
$ws.Range("D712").Value = 45075
$ws.Range("E712").Value = 13
$ws.Range("F712").Value = 100112030
$ws.Range("G712").Value = "Poroto granado"
$ws.Range("H712").Value = "Sin especificar"
$ws.Range("I712").Value = "Primera"
$ws.Range("J712").Value = 220
$ws.Range("K712").Value = 27000
$ws.Range("L712").Value = 30000
$ws.Range("M712").Value = 28636
$ws.Range("N712").Value = "`$/malla 25 kilos"
$ws.Range("O712").Value = "Provincia de Huasco"
$ws.Range("P712").Value = 1145
$ws.Range("Q712").Value = 25
$ws.Range("R712").Value = "Hortaliza"
